$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New overall_end (column L) values keyed by row number, per the 2029 roster update.
$overallEnd = @{
    2  = 86
    3  = 91
    5  = 91
    6  = 92
    8  = 83
    9  = 76
    12 = 86
    13 = 86
    16 = 81
    17 = 81
    23 = 80
    24 = 75
    26 = 93
    27 = 84
    29 = 86
    31 = 81
    32 = 88
    34 = 78
    36 = 79
    37 = 75
    39 = 88
    40 = 85
    41 = 76
    43 = 94
    45 = 80
    47 = 84
    48 = 85
    50 = 73
    52 = 91
    53 = 84
    58 = 94
    60 = 76
    61 = 93
    64 = 75
    67 = 92
    68 = 84
    71 = 89
    72 = 84
    73 = 84
    74 = 83
    76 = 77
    79 = 82
    80 = 82
    82 = 87
    84 = 69
    85 = 74
    86 = 76
}

foreach ($row in $overallEnd.Keys) {
    $ws.Cells.Item($row, 12).Value = $overallEnd[$row]
}

# Record the active/selected cell as L2, matching the author's saved selection.
$ws.Range("L2").Select()

$wb.Save()
